# Update the methods.summary sheet with refreshed TSS/AUC statistics and
# rotate the method labels in rows 7-9 (mda/rpart/glm -> rpart/glm/mda)
# per the backed-up run from 10-11 May 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (rf)
$ws.Range("C2").Value = 0.76366
$ws.Range("D2").Value = 0.09844368297374326
$ws.Range("E2").Value = 0.91752
$ws.Range("F2").Value = 0.04317273874497056

# Row 3 (brt)
$ws.Range("C3").Value = 0.702136
$ws.Range("D3").Value = 0.1024775722273866
$ws.Range("E3").Value = 0.884432
$ws.Range("F3").Value = 0.052885755949638555

# Row 4 (svm)
$ws.Range("C4").Value = 0.659018
$ws.Range("D4").Value = 0.09889014205871476
$ws.Range("E4").Value = 0.8604
$ws.Range("F4").Value = 0.05721215308762505

# Row 5 (gam)
$ws.Range("C5").Value = 0.64745
$ws.Range("D5").Value = 0.13654052192758204
$ws.Range("E5").Value = 0.82966
$ws.Range("F5").Value = 0.08167161628014731

# Row 6 (cart)
$ws.Range("C6").Value = 0.6448989898989899
$ws.Range("D6").Value = 0.1313711438937367
$ws.Range("E6").Value = 0.833939393939394
$ws.Range("F6").Value = 0.0785287959899357

# Row 7 (method label rotates mda -> rpart)
$ws.Range("B7").Value = "rpart"
$ws.Range("C7").Value = 0.603202
$ws.Range("D7").Value = 0.11489342101436778
$ws.Range("E7").Value = 0.830688
$ws.Range("F7").Value = 0.06932594629472662

# Row 8 (method label rotates rpart -> glm)
$ws.Range("B8").Value = "glm"
$ws.Range("C8").Value = 0.55956
$ws.Range("D8").Value = 0.10843905411296166
$ws.Range("E8").Value = 0.794496
$ws.Range("F8").Value = 0.06514674343634316

# Row 9 (method label rotates glm -> mda)
$ws.Range("B9").Value = "mda"
$ws.Range("C9").Value = 0.557744
$ws.Range("D9").Value = 0.11386531114507319
$ws.Range("E9").Value = 0.792884
$ws.Range("F9").Value = 0.07078737351037255

# Row 10 (fda)
$ws.Range("C10").Value = 0.502652
$ws.Range("D10").Value = 0.11640818109244085
$ws.Range("E10").Value = 0.762462
$ws.Range("F10").Value = 0.07569885858043543
